{"js": "// The document's first paragraph is an empty, bold, 12pt \"-Win---Researcher\"\n// paragraph used as top spacing. The edit duplicates that paragraph by\n// inserting an identical new empty paragraph immediately before it, so the\n// document ends up with two consecutive copies of that blank paragraph at\n// the very start of the body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Inserting an empty paragraph \"Before\" the existing first paragraph clones\n// its paragraph (pPr) and run (rPr) formatting, matching the diff exactly.\nfirstParagraph.insertParagraph(\"\", \"Before\");\n\nawait context.sync();\n", "ps1": "# The document's first paragraph is an empty, bold, 12pt \"-Win---Researcher\"\n# paragraph used as top spacing. The edit duplicates that paragraph by\n# inserting an identical new empty paragraph immediately before it, so the\n# document ends up with two consecutive copies of that blank paragraph at\n# the very start of the body.\n$d = $word.ActiveDocument\n$firstParagraph = $d.Paragraphs(1)\n\n# InsertParagraphBefore on the first paragraph's range clones its paragraph\n# (pPr) and run (rPr) formatting, matching the diff exactly.\n$firstParagraph.Range.InsertParagraphBefore()\n"}
